# Auto-generated edit script: updates FFXIV Chocobo Profits market-data cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 68.27
$ws.Cells.Item(15, 9).Value = 68.27
$ws.Cells.Item(15, 11).Value = 204.81
$ws.Cells.Item(15, 13).Value = -35.81

$ws.Cells.Item(39, 8).Value = 359
$ws.Cells.Item(39, 9).Value = 130.5
$ws.Cells.Item(39, 10).Value = 816
$ws.Cells.Item(39, 11).Value = 391.5
$ws.Cells.Item(39, 12).Value = 2448
$ws.Cells.Item(39, 13).Value = -95.5
$ws.Cells.Item(39, 14).Value = -3040

$ws.Cells.Item(43, 8).Value = 3126.5881
$ws.Cells.Item(43, 9).Value = 640.9091
$ws.Cells.Item(43, 10).Value = 7683.6665
$ws.Cells.Item(43, 11).Value = 640.9091
$ws.Cells.Item(43, 12).Value = 7683.6665
$ws.Cells.Item(43, 13).Value = -571.9091
$ws.Cells.Item(43, 14).Value = -7821.6665

$ws.Cells.Item(62, 8).Value = 1772.2858
$ws.Cells.Item(62, 9).Value = 1772.2858
$ws.Cells.Item(62, 11).Value = 1772.2858
$ws.Cells.Item(62, 13).Value = -1148.2858

$ws.Cells.Item(65, 8).Value = 1772.2858
$ws.Cells.Item(65, 9).Value = 1772.2858
$ws.Cells.Item(65, 11).Value = 8861.429
$ws.Cells.Item(65, 13).Value = -5741.429

$ws.Cells.Item(103, 8).Value = 6962
$ws.Cells.Item(103, 9).Value = 796.8
$ws.Cells.Item(103, 10).Value = 22375
$ws.Cells.Item(103, 11).Value = 2390.4
$ws.Cells.Item(103, 12).Value = 67125
$ws.Cells.Item(103, 13).Value = -1804.4
$ws.Cells.Item(103, 14).Value = -68297

$ws.Cells.Item(112, 8).Value = 1346.8182
$ws.Cells.Item(112, 9).Value = 800
$ws.Cells.Item(112, 10).Value = 1359.5349
$ws.Cells.Item(112, 11).Value = 2400
$ws.Cells.Item(112, 12).Value = 4078.6047
$ws.Cells.Item(112, 13).Value = -1292
$ws.Cells.Item(112, 14).Value = -6294.6047

$ws.Cells.Item(129, 8).Value = 1496.4828
$ws.Cells.Item(129, 10).Value = 1546.3455
$ws.Cells.Item(129, 12).Value = 4639.0365
$ws.Cells.Item(129, 14).Value = -14639.0365

$ws.Cells.Item(132, 8).Value = 24145126
$ws.Cells.Item(132, 9).Value = 32389180
$ws.Cells.Item(132, 11).Value = 97167540
$ws.Cells.Item(132, 13).Value = -97165010

$ws.Cells.Item(137, 8).Value = 664260.2
$ws.Cells.Item(137, 9).Value = 1539064.8
$ws.Cells.Item(137, 10).Value = 2822.5122
$ws.Cells.Item(137, 11).Value = 4617194.4
$ws.Cells.Item(137, 12).Value = 8467.536599999999
$ws.Cells.Item(137, 13).Value = -4614644.4
$ws.Cells.Item(137, 14).Value = -13567.5366

$ws.Cells.Item(141, 8).Value = 7894.4443
$ws.Cells.Item(141, 9).Value = 8640
$ws.Cells.Item(141, 10).Value = 4166.6665
$ws.Cells.Item(141, 11).Value = 25920
$ws.Cells.Item(141, 12).Value = 12499.9995
$ws.Cells.Item(141, 13).Value = -20740
$ws.Cells.Item(141, 14).Value = -22859.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(10, 8).Value = 17777
$ws.Cells.Item(10, 9).Value = 0
$ws.Cells.Item(10, 10).Value = 17777
$ws.Cells.Item(10, 11).Value = 0
$ws.Cells.Item(10, 12).Value = 17777
$ws.Cells.Item(10, 13).ClearContents()
$ws.Cells.Item(10, 14).Value = -18117

$ws.Cells.Item(32, 8).Value = 4591.279
$ws.Cells.Item(32, 9).Value = 4351.592
$ws.Cells.Item(32, 11).Value = 4351.592
$ws.Cells.Item(32, 13).Value = -4064.592

$ws.Cells.Item(132, 8).Value = 4106.8823
$ws.Cells.Item(132, 9).Value = 2591.3635
$ws.Cells.Item(132, 10).Value = 6885.3335
$ws.Cells.Item(132, 11).Value = 7774.0905
$ws.Cells.Item(132, 12).Value = 20656.0005
$ws.Cells.Item(132, 13).Value = -5244.0905
$ws.Cells.Item(132, 14).Value = -25716.0005

$ws.Cells.Item(137, 8).Value = 40704
$ws.Cells.Item(137, 10).Value = 40704
$ws.Cells.Item(137, 12).Value = 40704
$ws.Cells.Item(137, 14).Value = -50904

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 1093.7646
$ws.Cells.Item(107, 9).Value = 1009.7308
$ws.Cells.Item(107, 10).Value = 1366.875
$ws.Cells.Item(107, 11).Value = 1009.7308
$ws.Cells.Item(107, 12).Value = 1366.875
$ws.Cells.Item(107, 13).Value = 910.2692
$ws.Cells.Item(107, 14).Value = -5206.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 2652
$ws.Cells.Item(2, 9).Value = 2652
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 2652
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = -2539
$ws.Cells.Item(2, 14).ClearContents()

$ws.Cells.Item(31, 8).Value = 194862.3
$ws.Cells.Item(31, 9).Value = 466819.8
$ws.Cells.Item(31, 10).Value = 2502.0977
$ws.Cells.Item(31, 11).Value = 466819.8
$ws.Cells.Item(31, 12).Value = 2502.0977
$ws.Cells.Item(31, 13).Value = -466524.8
$ws.Cells.Item(31, 14).Value = -3092.0977

$ws.Cells.Item(34, 8).Value = 194862.3
$ws.Cells.Item(34, 9).Value = 466819.8
$ws.Cells.Item(34, 10).Value = 2502.0977
$ws.Cells.Item(34, 11).Value = 466819.8
$ws.Cells.Item(34, 12).Value = 2502.0977
$ws.Cells.Item(34, 13).Value = -466617.8
$ws.Cells.Item(34, 14).Value = -2906.0977

$ws.Cells.Item(137, 8).Value = 44894.285
$ws.Cells.Item(137, 10).Value = 44894.285
$ws.Cells.Item(137, 12).Value = 44894.285
$ws.Cells.Item(137, 14).Value = -55094.285

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 1952
$ws.Cells.Item(22, 10).Value = 1952
$ws.Cells.Item(22, 12).Value = 5856
$ws.Cells.Item(22, 14).Value = -6194

$ws.Cells.Item(27, 8).Value = 1952
$ws.Cells.Item(27, 10).Value = 1952
$ws.Cells.Item(27, 12).Value = 5856
$ws.Cells.Item(27, 14).Value = -6060

$ws.Cells.Item(55, 8).Value = 9395.714
$ws.Cells.Item(55, 9).Value = 0
$ws.Cells.Item(55, 10).Value = 9395.714
$ws.Cells.Item(55, 11).Value = 0
$ws.Cells.Item(55, 12).Value = 28187.142
$ws.Cells.Item(55, 13).ClearContents()
$ws.Cells.Item(55, 14).Value = -28541.142

$ws.Cells.Item(68, 8).Value = 2650.3242
$ws.Cells.Item(68, 9).Value = 878.24
$ws.Cells.Item(68, 10).Value = 3554.449
$ws.Cells.Item(68, 11).Value = 2634.72
$ws.Cells.Item(68, 12).Value = 10663.347
$ws.Cells.Item(68, 13).Value = -1823.72
$ws.Cells.Item(68, 14).Value = -12285.347

$ws.Cells.Item(71, 8).Value = 2650.3242
$ws.Cells.Item(71, 9).Value = 878.24
$ws.Cells.Item(71, 10).Value = 3554.449
$ws.Cells.Item(71, 11).Value = 7904.16
$ws.Cells.Item(71, 12).Value = 31990.041
$ws.Cells.Item(71, 13).Value = -3848.16
$ws.Cells.Item(71, 14).Value = -40102.041

$ws.Cells.Item(107, 8).Value = 6037137
$ws.Cells.Item(107, 9).Value = 415.21054
$ws.Cells.Item(107, 10).Value = 11134814
$ws.Cells.Item(107, 11).Value = 1245.63162
$ws.Cells.Item(107, 12).Value = 33404442
$ws.Cells.Item(107, 13).Value = 674.3683800000001
$ws.Cells.Item(107, 14).Value = -33408282

$ws.Cells.Item(122, 8).Value = 2179.4238
$ws.Cells.Item(122, 10).Value = 2867.45
$ws.Cells.Item(122, 12).Value = 25807.05
$ws.Cells.Item(122, 14).Value = -30707.05

$ws.Cells.Item(131, 8).Value = 886.99
$ws.Cells.Item(131, 9).Value = 1261.3334
$ws.Cells.Item(131, 10).Value = 875.4123499999999
$ws.Cells.Item(131, 11).Value = 3784.0002
$ws.Cells.Item(131, 12).Value = 2626.23705
$ws.Cells.Item(131, 13).Value = 1255.9998
$ws.Cells.Item(131, 14).Value = -12706.23705

$ws.Cells.Item(132, 8).Value = 2217.85
$ws.Cells.Item(132, 9).Value = 940.8
$ws.Cells.Item(132, 11).Value = 8467.199999999999
$ws.Cells.Item(132, 13).Value = -5937.199999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(7, 8).Value = 1112847.5
$ws.Cells.Item(7, 10).Value = 5007499
$ws.Cells.Item(7, 12).Value = 5007499
$ws.Cells.Item(7, 14).Value = -5007723

$ws.Cells.Item(8, 8).Value = 1112847.5
$ws.Cells.Item(8, 10).Value = 5007499
$ws.Cells.Item(8, 12).Value = 5007499
$ws.Cells.Item(8, 14).Value = -5007777

$ws.Cells.Item(12, 8).Value = 24199.666
$ws.Cells.Item(12, 10).Value = 29799.5
$ws.Cells.Item(12, 12).Value = 29799.5
$ws.Cells.Item(12, 14).Value = -30079.5

$ws.Cells.Item(132, 8).Value = 2170.3274
$ws.Cells.Item(132, 9).Value = 1085.1875
$ws.Cells.Item(132, 10).Value = 3680.087
$ws.Cells.Item(132, 11).Value = 3255.5625
$ws.Cells.Item(132, 12).Value = 11040.261
$ws.Cells.Item(132, 13).Value = -725.5625
$ws.Cells.Item(132, 14).Value = -16100.261

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(13, 8).Value = 18000
$ws.Cells.Item(13, 10).Value = 18000
$ws.Cells.Item(13, 12).Value = 18000
$ws.Cells.Item(13, 14).Value = -18280

$ws.Cells.Item(122, 8).Value = 4791.625
$ws.Cells.Item(122, 9).Value = 2706
$ws.Cells.Item(122, 11).Value = 8118
$ws.Cells.Item(122, 13).Value = -5668

$ws.Cells.Item(132, 8).Value = 3761.2354
$ws.Cells.Item(132, 9).Value = 2993.25
$ws.Cells.Item(132, 10).Value = 4858.357
$ws.Cells.Item(132, 11).Value = 8979.75
$ws.Cells.Item(132, 12).Value = 14575.071
$ws.Cells.Item(132, 13).Value = -6449.75
$ws.Cells.Item(132, 14).Value = -19635.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(8, 8).Value = 335668.66
$ws.Cells.Item(8, 9).Value = 335668.66
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 335668.66
$ws.Cells.Item(8, 12).Value = 0
$ws.Cells.Item(8, 13).Value = -335528.66
$ws.Cells.Item(8, 14).ClearContents()

$ws.Cells.Item(11, 8).Value = 762500
$ws.Cells.Item(11, 9).Value = 1000000
$ws.Cells.Item(11, 11).Value = 1000000
$ws.Cells.Item(11, 13).Value = -999858
